$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "56.065.46"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  +3.32%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "2.316.96"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  +2.27%  "
$ws.Range("E4").Value = "  +0.24%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "516.70"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  +3.91%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "132.93"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +3.26%  "
$ws.Range("E7").Value = "  -0.27%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.535"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "  +1.95%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "2.340.96"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +3.03%  "
$ws.Range("E10").Value = "  +8.24%  "
$ws.Range("E11").Value = "  +0.92%  "
$ws.Range("E12").Value = "  +7.96%  "
$ws.Range("E13").Value = "  +1.78%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "24.05"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  +5.19%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "2.732.91"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  +2.46%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "56.259.34"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  +3.75%  "
$ws.Range("E17").Value = "  +4.43%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "2.329.48"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  +2.54%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "10.54"
$ws.Range("D19").Style = "Normal"
$ws.Range("E20").Value = "  +2.80%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "322.21"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  +6.36%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "6.68"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +5.89%  "
$ws.Range("E23").Value = "  +0.02%  "
$ws.Range("D24").NumberFormat = "@"
$ws.Range("D24").Value = "60.80"
$ws.Range("D24").Style = "Normal"
$ws.Range("E24").Value = "  -0.25%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "0.994"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "  -0.28%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "0.158"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +5.79%  "
$ws.Range("E27").Value = "  +4.05%  "
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "172.07"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +0.62%  "
$ws.Range("E29").Value = "  +9.82%  "
$ws.Range("B30").Value = "PEPE"
$ws.Range("C30").Value = "https://coinranking.com/coin/03WI8NQPF+pepe-pepe"
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "0.0₃0723"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  +4.65%  "
$ws.Range("B31").Value = "PancakeSwap"
$ws.Range("C31").Value = "https://coinranking.com/coin/ncYFcP709+pancakeswap-cake"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "1.68"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  +4.57%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "6.23"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +5.43%  "
$ws.Range("E33").Value = "  +3.37%  "
$ws.Range("E34").Value = "  -0.02%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "0.992"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -0.45%  "
$ws.Range("E36").Value = "  +5.77%  "
$ws.Range("E37").Value = "  -1.19%  "
$ws.Range("E38").Value = "  +7.70%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "37.36"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "  +3.87%  "
$ws.Range("E41").Value = "  +2.50%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "3.61"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "  +7.46%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "138.81"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +10.59%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "5.12"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "  +6.86%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "268.30"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +10.96%  "
$ws.Range("E46").Value = "  +2.96%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0925"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  +3.74%  "
$ws.Range("E48").Value = "  +1.57%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.382"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  +2.53%  "
$ws.Range("E50").Value = "  +5.09%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "16.94"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  +4.96%  "
